$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B22").Value = "<200 Manuf., Constr., and Mining, `n<100 Others"
$ws.Range("C23").Value = "<NT`$80Millionlion Manuf., Constr., and Mining, `n<NT`$100Millionlion Others"
$ws.Range("B24").Value = ">200 Mnf., CnsTurnover, & Minin., `n>100 Others"
$ws.Range("C24").Value = ">=NT`$80Millionlion Manuf., Constr., and Mining,`n >=NT`$100Millionlion Others"
